$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 33 with the same structure as the existing data rows.
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 33
$ws.Range("C33").Value = 11
$ws.Range("D33").Value = 13
$ws.Range("E33").Value = 37
$ws.Range("F33").Value = 57
$ws.Range("G33").Value = 94
